# Apply the "new book edition" content edits to the first body paragraph:
#   "A plain paragraph with some bold and some italic"
# becomes
#   "A plain paragraph with some bold text and some italic"
# and a "_GoBack" bookmark (Word's standard "last edit location" marker)
# ends up wrapped around the "italic" run.

$d = $word.ActiveDocument

# 1) Merge "A plain paragraph with" + " some " into a single run by
#    replacing across the run boundary with identical text.
$d.Content.Find.Execute("with some ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "with some ", 2)

# 2) Insert the word "text" between "bold" and "and some ".
$d.Content.Find.Execute(" and some ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " text and some ", 2)

# 3) Re-establish Word's "_GoBack" bookmark around the last edited text
#    (the "italic" run), mirroring what Word stamps on save after an edit.
$r = $d.Content
$r.Find.Execute("italic", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $r)
